$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update end date of "Detailed Design" task (row 7) from 11/16/2012 to 11/10/2012
$ws.Range("E7").Value = 41223

# Insert a new row at row 15 for the new task "SE_N15_DS_06", shifting
# the existing rows 15-26 down to 16-27
$ws.Rows("15").Insert()

# Populate the newly inserted row 15 with the new task details
$ws.Range("A15").Value = "SE_N15_DS_06"
$ws.Range("C15").Value = "Kamil Mrowic"
$ws.Range("D15").Value = 41214
$ws.Range("E15").Value = 41221
$ws.Range("F15").Value = "Researching DBMS Systems (Prosgres vs MySql)"
